# Update automatico via Actualizar 03-07-2021 13-15-43
# This mirrors a scheduled refresh that re-stamps the "last updated" timestamp
# (column D) for each block of rows, shifting the previous values down one
# generation and writing the brand-new timestamp into the most recent block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newest = 44262.55235347203
$shift1 = 44262.53097673611
$shift2 = 44262.50961724537

for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = $newest
}

for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = $shift1
}

for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = $shift2
}
